$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two string values: A3 ("Saturday's Market", curly apostrophe) and A4 ("Krux")
# so that A3 becomes "Krux" and A4 becomes "Saturday's Market" with a straight apostrophe.
$ws.Range("A3").Value = "Krux"
$ws.Range("A4").Value = "Saturday's Market"

# Update the selected/active cell from B5 to A2
$ws.Range("A2").Select()
